$d = $word.ActiveDocument

# 1. "Javascript, JQuery, NODE JS, MongoDB and SEO tools." -> "Javascript, CMS and SEO tools."
$d.Content.Find.Execute("JQuery, NODE JS, MongoDB and SEO tools.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CMS and SEO tools.", 2)

# 2. "HTML, CSS, PHP." -> "HTML, CSS, Javascript and PHP."
$d.Content.Find.Execute(", PHP.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", Javascript and PHP.", 2)

# 3. Typo fix: "tchniques" -> "techniques"
$d.Content.Find.Execute("tchniques", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "techniques", 2)
